$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column values (token text) where changed
$ws.Range("B3").Value  = "<number>"
$ws.Range("B9").Value  = "<its>"
$ws.Range("B13").Value = "<six>"
$ws.Range("B14").Value = "<alt>"
$ws.Range("B16").Value = "<yankee>"
$ws.Range("B17").Value = "<sentence>"
$ws.Range("B18").Value = "<tab>"

# Update C column values (counts)
$ws.Range("C2").Value  = 24
$ws.Range("C3").Value  = 26
$ws.Range("C4").Value  = 21
$ws.Range("C5").Value  = 23
$ws.Range("C6").Value  = 16
$ws.Range("C7").Value  = 14
$ws.Range("C8").Value  = 20
$ws.Range("C9").Value  = 27
$ws.Range("C10").Value = 27
$ws.Range("C11").Value = 26
$ws.Range("C12").Value = 29
$ws.Range("C13").Value = 31
$ws.Range("C14").Value = 19
$ws.Range("C15").Value = 20
$ws.Range("C16").Value = 26
$ws.Range("C17").Value = 22
$ws.Range("C18").Value = 22
